$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.095.07'
$ws.Range("E2").Value = '  +1.18%  '
$ws.Range("D3").Value = '2.218.82'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'290.60"
$ws.Range("E5").Value = '  -2.70%  '
$ws.Range("D6").Value = "'88.44"
$ws.Range("E6").Value = '  +5.59%  '
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  +1.34%  '
$ws.Range("D10").Value = "'30.86"
$ws.Range("E10").Value = '  +4.13%  '
$ws.Range("D11").Value = "'0.0783"
$ws.Range("E11").Value = '  +0.30%  '
$ws.Range("E12").Value = '  +3.61%  '
$ws.Range("E13").Value = '  +2.62%  '
$ws.Range("D14").Value = "'6.48"
$ws.Range("E14").Value = '  +2.89%  '
$ws.Range("D15").Value = '2.560.00'
$ws.Range("E15").Value = '  +0.09%  '
$ws.Range("D16").Value = "'14.05"
$ws.Range("E16").Value = '  -0.50%  '
$ws.Range("D17").Value = '2.230.41'
$ws.Range("E17").Value = '  +0.73%  '
$ws.Range("D18").Value = "'0.730"
$ws.Range("E18").Value = '  +1.64%  '
$ws.Range("D19").Value = '40.027.95'
$ws.Range("E19").Value = '  +1.23%  '
$ws.Range("D20").Value = "'11.92"
$ws.Range("E20").Value = '  +14.47%  '
$ws.Range("E21").Value = '  +0.93%  '
$ws.Range("D22").Value = "'5.83"
$ws.Range("E22").Value = '  +1.48%  '
$ws.Range("D23").Value = "'65.66"
$ws.Range("E23").Value = '  +0.95%  '
$ws.Range("D24").Value = "'235.82"
$ws.Range("E24").Value = '  +1.57%  '
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("D26").Value = "'2.47"
$ws.Range("E26").Value = '  +1.80%  '
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("D28").Value = "'22.67"
$ws.Range("E28").Value = '  -0.31%  '
$ws.Range("E29").Value = '  +0.99%  '
$ws.Range("D30").Value = "'9.26"
$ws.Range("E30").Value = '  +1.02%  '
$ws.Range("D31").Value = "'152.99"
$ws.Range("E31").Value = '  +2.44%  '
$ws.Range("D32").Value = "'32.25"
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = '  -0.07%  '
$ws.Range("D34").Value = "'4.97"
$ws.Range("E34").Value = '  +2.70%  '
$ws.Range("D35").Value = "'0.0722"
$ws.Range("E35").Value = '  +2.97%  '
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("E37").Value = '  +6.97%  '
$ws.Range("E38").Value = '  +0.61%  '
$ws.Range("E39").Value = '  -0.96%  '
$ws.Range("D40").Value = "'0.1000"
$ws.Range("E40").Value = '  +2.99%  '
$ws.Range("E41").Value = '  +3.36%  '
$ws.Range("D42").Value = '2.097.89'
$ws.Range("E42").Value = '  +8.85%  '
$ws.Range("E43").Value = '  +4.93%  '
$ws.Range("E44").Value = '  +1.76%  '
$ws.Range("E45").Value = '  +1.12%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = "'17.72"
$ws.Range("E46").Value = '  +8.98%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = "'9.87"
$ws.Range("E47").Value = '  +7.07%  '
$ws.Range("E48").Value = '  +1.85%  '
$ws.Range("D49").Value = '2.431.98'
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("D50").Value = "'69.65"
$ws.Range("E50").Value = '  -1.79%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = "'88.94"
$ws.Range("E51").Value = '  +0.25%  '
